# Generate Report for Archive
# Update the Status value from "Ready for handoff" to "In Translation"
# wherever it appears, and shrink the corresponding Status column widths
# to match the narrower text.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value2 = "In Translation"
        }
    }
}

# Narrow the "Status" columns that previously held the wider text.
# (12.5 is the input that lands closest on this engine's column-width grid
# to the target stored width of 13.4101848602295.)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
